$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 'sd'
$ws.Range("J2").Value = 'Statement-non-opinion'
$ws.Range("I7").Value = 'sv'
$ws.Range("J7").Value = 'Statement-opinion'
$ws.Range("I12").Value = 'b'
$ws.Range("J12").Value = 'Acknowledge (Backchannel)'
$ws.Range("I13").Value = 'sd'
$ws.Range("J13").Value = 'Statement-non-opinion'
$ws.Range("I21").Value = '%'
$ws.Range("J21").Value = 'Uninterpretable'
$ws.Range("I22").Value = 'ba'
$ws.Range("J22").Value = 'Appreciation'
$ws.Range("I23").Value = 'sd'
$ws.Range("J23").Value = 'Statement-non-opinion'
$ws.Range("I31").Value = 'b'
$ws.Range("J31").Value = 'Acknowledge (Backchannel)'
$ws.Range("I33").Value = 'sd'
$ws.Range("J33").Value = 'Statement-non-opinion'
$ws.Range("I34").Value = 'sd'
$ws.Range("J34").Value = 'Statement-non-opinion'
$ws.Range("I38").Value = 'ba'
$ws.Range("J38").Value = 'Appreciation'
$ws.Range("I45").Value = 'ba'
$ws.Range("J45").Value = 'Appreciation'
$ws.Range("I58").Value = 'aa'
$ws.Range("J58").Value = 'Agree/Accept'
$ws.Range("I60").Value = 'b'
$ws.Range("J60").Value = 'Acknowledge (Backchannel)'
$ws.Range("I61").Value = 'sd'
$ws.Range("J61").Value = 'Statement-non-opinion'
$ws.Range("I64").Value = 'qy'
$ws.Range("J64").Value = 'Yes-No-Question'
$ws.Range("I78").Value = 'b'
$ws.Range("J78").Value = 'Acknowledge (Backchannel)'
$ws.Range("I82").Value = 'b'
$ws.Range("J82").Value = 'Acknowledge (Backchannel)'
$ws.Range("I125").Value = 'aa'
$ws.Range("J125").Value = 'Agree/Accept'
$ws.Range("I128").Value = '%'
$ws.Range("J128").Value = 'Uninterpretable'
$ws.Range("I130").Value = 'sd'
$ws.Range("J130").Value = 'Statement-non-opinion'
$ws.Range("I132").Value = 'sd'
$ws.Range("J132").Value = 'Statement-non-opinion'
$ws.Range("I133").Value = 'sv'
$ws.Range("J133").Value = 'Statement-opinion'
$ws.Range("I134").Value = 'ba'
$ws.Range("J134").Value = 'Appreciation'
$ws.Range("I139").Value = 'sv'
$ws.Range("J139").Value = 'Statement-opinion'
$ws.Range("I141").Value = '%'
$ws.Range("J141").Value = 'Uninterpretable'
$ws.Range("I151").Value = 'sd'
$ws.Range("J151").Value = 'Statement-non-opinion'
$ws.Range("I152").Value = '%'
$ws.Range("J152").Value = 'Uninterpretable'
$ws.Range("I153").Value = 'sd'
$ws.Range("J153").Value = 'Statement-non-opinion'
$ws.Range("I154").Value = 'b'
$ws.Range("J154").Value = 'Acknowledge (Backchannel)'
$ws.Range("I156").Value = 'b'
$ws.Range("J156").Value = 'Acknowledge (Backchannel)'
$ws.Range("I164").Value = 'ba'
$ws.Range("J164").Value = 'Appreciation'
$ws.Range("I165").Value = 'aa'
$ws.Range("J165").Value = 'Agree/Accept'
$ws.Range("I166").Value = 'aa'
$ws.Range("J166").Value = 'Agree/Accept'
$ws.Range("I170").Value = 'aa'
$ws.Range("J170").Value = 'Agree/Accept'
$ws.Range("I176").Value = 'aa'
$ws.Range("J176").Value = 'Agree/Accept'
$ws.Range("I180").Value = 'sv'
$ws.Range("J180").Value = 'Statement-opinion'
$ws.Range("I195").Value = 'b'
$ws.Range("J195").Value = 'Acknowledge (Backchannel)'
$ws.Range("I196").Value = '%'
$ws.Range("J196").Value = 'Uninterpretable'
$ws.Range("I213").Value = 'sv'
$ws.Range("J213").Value = 'Statement-opinion'
$ws.Range("I232").Value = 'sd'
$ws.Range("J232").Value = 'Statement-non-opinion'
$ws.Range("I253").Value = 'sd'
$ws.Range("J253").Value = 'Statement-non-opinion'
$ws.Range("I269").Value = 'sd'
$ws.Range("J269").Value = 'Statement-non-opinion'
$ws.Range("I274").Value = 'sv'
$ws.Range("J274").Value = 'Statement-opinion'
$ws.Range("I300").Value = 'sv'
$ws.Range("J300").Value = 'Statement-opinion'
$ws.Range("I311").Value = 'sd'
$ws.Range("J311").Value = 'Statement-non-opinion'
$ws.Range("I316").Value = 'aa'
$ws.Range("J316").Value = 'Agree/Accept'
$ws.Range("I317").Value = 'sd'
$ws.Range("J317").Value = 'Statement-non-opinion'
$ws.Range("I323").Value = 'sv'
$ws.Range("J323").Value = 'Statement-opinion'
$ws.Range("I360").Value = 'sv'
$ws.Range("J360").Value = 'Statement-opinion'
$ws.Range("I364").Value = 'aa'
$ws.Range("J364").Value = 'Agree/Accept'
$ws.Range("I373").Value = 'sd'
$ws.Range("J373").Value = 'Statement-non-opinion'
